$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (want-to-go count) in column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 481
$ws1.Range("F7").Value  = 838
$ws1.Range("F10").Value = 2088
$ws1.Range("F14").Value = 982
$ws1.Range("F18").Value = 10426
$ws1.Range("F19").Value = 1042

# Sheet "本地生活" (local life) - update same counter
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5658
$ws3.Range("F3").Value = 459

# Sheet "全部类型" (all types, aggregated view of all sheets)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 5658
$ws4.Range("F4").Value  = 459
$ws4.Range("F11").Value = 481
$ws4.Range("F12").Value = 838
$ws4.Range("F16").Value = 2088
$ws4.Range("F22").Value = 982
$ws4.Range("F31").Value = 1042
